# Auto-generated COM-interop script applying the scheduled-runner price refresh
# across the Aegis_Profits workbook (see commit message: "chore: update Sheets via scheduled runner").
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 713
$ws.Range("I28").Value = 290.07693
$ws.Range("J28").Value = 1105.7142
$ws.Range("K28").Value = 290.07693
$ws.Range("L28").Value = 1105.7142
$ws.Range("M28").Value = 194.92307
$ws.Range("N28").Value = -2075.7142
$ws.Range("H137").Value = 1340.907
$ws.Range("I137").Value = 1376.0646
$ws.Range("J137").Value = 1250.0834
$ws.Range("K137").Value = 4128.1938
$ws.Range("L137").Value = 3750.2502
$ws.Range("M137").Value = -1578.1938
$ws.Range("N137").Value = -8850.2502

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22425.938
$ws.Range("I32").Value = 3826.4028
$ws.Range("K32").Value = 3826.4028
$ws.Range("M32").Value = -3539.4028
$ws.Range("H61").Value = 1027.5807
$ws.Range("I61").Value = 905.46155
$ws.Range("J61").Value = 1662.6
$ws.Range("K61").Value = 905.46155
$ws.Range("L61").Value = 1662.6
$ws.Range("M61").Value = -693.46155
$ws.Range("N61").Value = -2086.6
$ws.Range("H74").Value = 2361.3125
$ws.Range("I74").Value = 1437.0869
$ws.Range("J74").Value = 4723.222
$ws.Range("K74").Value = 1437.0869
$ws.Range("L74").Value = 4723.222
$ws.Range("M74").Value = -563.0869
$ws.Range("N74").Value = -6471.222
$ws.Range("H77").Value = 2361.3125
$ws.Range("I77").Value = 1437.0869
$ws.Range("J77").Value = 4723.222
$ws.Range("K77").Value = 7185.4345
$ws.Range("L77").Value = 23616.11
$ws.Range("M77").Value = -2817.4345
$ws.Range("N77").Value = -32352.11
$ws.Range("H97").Value = 37826.258
$ws.Range("I97").Value = 39263.81
$ws.Range("J97").Value = 450
$ws.Range("K97").Value = 39263.81
$ws.Range("L97").Value = 450
$ws.Range("M97").Value = -38767.81
$ws.Range("N97").Value = -1442
$ws.Range("H98").Value = 6200
$ws.Range("J98").Value = 6200
$ws.Range("L98").Value = 6200
$ws.Range("N98").Value = -12190
$ws.Range("H132").Value = 3337.4517
$ws.Range("I132").Value = 3272.7036
$ws.Range("K132").Value = 9818.110799999999
$ws.Range("M132").Value = -7288.110799999999
$ws.Range("H136").Value = 1027.5807
$ws.Range("I136").Value = 905.46155
$ws.Range("J136").Value = 1662.6
$ws.Range("K136").Value = 2716.38465
$ws.Range("L136").Value = 4987.799999999999
$ws.Range("M136").Value = -166.38465
$ws.Range("N136").Value = -10087.8

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2403.5667
$ws.Range("I134").Value = 2408.6428
$ws.Range("J134").Value = 2332.5
$ws.Range("K134").Value = 7225.928400000001
$ws.Range("L134").Value = 6997.5
$ws.Range("M134").Value = -4690.928400000001
$ws.Range("N134").Value = -12067.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 21618.744
$ws.Range("I31").Value = 30751.295
$ws.Range("K31").Value = 30751.295
$ws.Range("M31").Value = -30456.295
$ws.Range("H34").Value = 21618.744
$ws.Range("I34").Value = 30751.295
$ws.Range("K34").Value = 30751.295
$ws.Range("M34").Value = -30549.295
$ws.Range("H45").Value = 13333.333
$ws.Range("H58").Value = 7731.6313
$ws.Range("I58").Value = 1250.0769
$ws.Range("J58").Value = 21775
$ws.Range("K58").Value = 1250.0769
$ws.Range("L58").Value = 21775
$ws.Range("M58").Value = -1047.0769
$ws.Range("N58").Value = -22181
$ws.Range("H62").Value = 2590.4167
$ws.Range("I62").Value = 2498.3333
$ws.Range("K62").Value = 2498.3333
$ws.Range("M62").Value = -1874.3333
$ws.Range("H65").Value = 2590.4167
$ws.Range("I65").Value = 2498.3333
$ws.Range("K65").Value = 12491.6665
$ws.Range("M65").Value = -9371.666499999999
$ws.Range("H94").Value = 1286.6666
$ws.Range("I94").Value = 1025
$ws.Range("J94").Value = 1361.4286
$ws.Range("K94").Value = 1025
$ws.Range("L94").Value = 1361.4286
$ws.Range("M94").Value = -574
$ws.Range("N94").Value = -2263.4286
$ws.Range("H132").Value = 51727296
$ws.Range("I132").Value = 45457540
$ws.Range("J132").Value = 71432250
$ws.Range("K132").Value = 136372620
$ws.Range("L132").Value = 214296750
$ws.Range("M132").Value = -136370090
$ws.Range("N132").Value = -214301810
$ws.Range("H134").Value = 1692.9524
$ws.Range("I134").Value = 1869
$ws.Range("J134").Value = 1406.875
$ws.Range("K134").Value = 5607
$ws.Range("L134").Value = 4220.625
$ws.Range("M134").Value = -3072
$ws.Range("N134").Value = -9290.625
$ws.Range("H136").Value = 7731.6313
$ws.Range("I136").Value = 1250.0769
$ws.Range("J136").Value = 21775
$ws.Range("K136").Value = 3750.2307
$ws.Range("L136").Value = 65325
$ws.Range("M136").Value = -1200.2307
$ws.Range("N136").Value = -70425

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 3603.3333
$ws.Range("I123").Value = 2222.5
$ws.Range("J123").Value = 4708
$ws.Range("K123").Value = 6667.5
$ws.Range("L123").Value = 14124
$ws.Range("M123").Value = -4217.5
$ws.Range("N123").Value = -19024

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 724758.5600000001
$ws.Range("I5").Value = 5000000
$ws.Range("J5").Value = 12218.333
$ws.Range("K5").Value = 5000000
$ws.Range("L5").Value = 12218.333
$ws.Range("M5").Value = -4999888
$ws.Range("N5").Value = -12442.333
$ws.Range("H100").Value = 42980
$ws.Range("J100").Value = 42980
$ws.Range("L100").Value = 42980
$ws.Range("N100").Value = -45144
$ws.Range("H122").Value = 1333.3334
$ws.Range("I122").Value = 1333.3334
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4000.0002
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1550.0002
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 2312.2593
$ws.Range("I132").Value = 1684.6666
$ws.Range("J132").Value = 3567.4443
$ws.Range("K132").Value = 5053.9998
$ws.Range("L132").Value = 10702.3329
$ws.Range("M132").Value = -2523.9998
$ws.Range("N132").Value = -15762.3329

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2552.9534
$ws.Range("I132").Value = 2750.8965
$ws.Range("J132").Value = 2142.9285
$ws.Range("K132").Value = 8252.6895
$ws.Range("L132").Value = 6428.7855
$ws.Range("M132").Value = -5722.6895
$ws.Range("N132").Value = -11488.7855
$ws.Range("H136").Value = 1241.9667
$ws.Range("I136").Value = 1191.3928
$ws.Range("J136").Value = 1950
$ws.Range("K136").Value = 3574.1784
$ws.Range("L136").Value = 5850
$ws.Range("M136").Value = -1024.1784
$ws.Range("N136").Value = -10950

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2998.5652
$ws.Range("I132").Value = 2937.842
$ws.Range("J132").Value = 3287
$ws.Range("K132").Value = 8813.526
$ws.Range("L132").Value = 9861
$ws.Range("M132").Value = -6283.526
$ws.Range("N132").Value = -14921
$ws.Range("H136").Value = 851.53845
$ws.Range("I136").Value = 724.7222
$ws.Range("J136").Value = 1136.875
$ws.Range("K136").Value = 2174.1666
$ws.Range("L136").Value = 3410.625
$ws.Range("M136").Value = 375.8334
$ws.Range("N136").Value = -8510.625
